# 25-09-19 2a Visualización corregida y completa.
# Update "estatus" (column O) values on the active sheet for several rows,
# correcting the previously mis-assigned status labels.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("O3").Value  = "Contacto"
$ws.Range("O9").Value  = "Afiliado"
$ws.Range("O12").Value = "Afiliado"
$ws.Range("O14").Value = "Contacto"
$ws.Range("O15").Value = "Afiliado"
$ws.Range("O17").Value = "Afiliado"

# Move the sheet's active cell/selection to O3, matching the cursor
# position left behind after the correction.
[void]$ws.Range("O3").Select()
